# Updated symbol list (coin prices / ranking reshuffle) to match the
# Thu Dec 22 23:29:17 UTC 2022 GitHub Actions run.
#
# Many of the new values (e.g. "245.77", "0.0001500", "0.00000000750")
# look like numbers but must stay as plain text, exactly like the rest
# of the "Price" column in this sheet. Assigning them directly would make
# Excel auto-convert them to numeric cells (losing trailing zeros / exact
# formatting), so we prefix those with a text-quote ('), then restore the
# cell to the "Normal" style so no stray NumberFormat/style gets attached.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (BNB) ---
$ws.Range("D2").Value = "'245.77"

# --- Row 3 (OKB) ---
$ws.Range("D3").Value = "'21.58"

# --- Row 4 (HuobiToken) ---
$ws.Range("D4").Value = "'5.439"

# --- Row 5 (Cronos) ---
$ws.Range("D5").Value = "'0.05769"

# --- Row 7 (KuCoinToken) ---
$ws.Range("D7").Value = "'6.347"

# --- Row 8 (MXToken) ---
$ws.Range("D8").Value = "'0.8189"

# --- Row 9 (FTXToken) ---
$ws.Range("D9").Value = "'1.007"
$ws.Range("E9").Value = "8FTXTokenFTT"

# --- Row 10 (was WazirX -> now One) ---
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "'0.01104"
$ws.Range("E10").Value = "9OneONEBestin24h"

# --- Row 11 (was MandalaExchangeToken -> now WazirX) ---
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1429"
$ws.Range("E11").Value = "10WazirXWRX"

# --- Row 12 (was LiechtensteinCryptoassetsExchange -> now MandalaExchangeToken) ---
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.07276"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"

# --- Row 13 (was ProBitToken -> now LiechtensteinCryptoassetsExchange) ---
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "'0.03109"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"

# --- Row 14 (BitrueCoin) ---
$ws.Range("D14").Value = "'0.03110"

# --- Row 15 (MCDex) ---
$ws.Range("D15").Value = "'4.162"

# --- Row 16 (BitMartToken) ---
$ws.Range("D16").Value = "'0.09396"

# --- Row 17 (BitForexToken) ---
$ws.Range("D17").Value = "'0.001590"

# --- Row 18 (CoinExToken) ---
$ws.Range("D18").Value = "'0.04804"

# --- Row 19 (was One -> now TigerCash) ---
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D19").Value = "'0.006272"
$ws.Range("E19").Value = "18TigerCashTCH"

# --- Row 20 (was TigerCash -> now HotbitToken) ---
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").Value = "'0.004121"
$ws.Range("E20").Value = "19HotbitTokenHTB"

# --- Row 21 (was HotbitToken -> now BitKan) ---
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D21").Value = "'0.0009938"
$ws.Range("E21").Value = "20BitKanKAN"

# --- Row 22 (was BitKan -> now NitroEx) ---
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").Value = "'0.0001500"
$ws.Range("E22").Value = "21NitroExNTX"

# --- Row 23 (was NitroEx -> now LEO) ---
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").Value = "'3.743"
$ws.Range("E23").Value = "22LEOLEO"

# --- Row 24 (was LEO -> now BTSEToken) ---
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").Value = "'2.200"
$ws.Range("E24").Value = "23BTSETokenBTSE"

# --- Row 25 (was BTSEToken -> now BitpandaEcosystemToken) ---
$ws.Range("B25").Value = "BitpandaEcosystemToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D25").Value = "'0.3181"
$ws.Range("E25").Value = "24BitpandaEcosystemTokenBEST"

# --- Row 26 (was BitpandaEcosystemToken -> now ProBitToken) ---
$ws.Range("B26").Value = "ProBitToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D26").Value = "'0.1329"
$ws.Range("E26").Value = "25ProBitTokenPROB"

# --- Row 27 (UpBots) ---
$ws.Range("D27").Value = "'0.0003998"

# --- Row 40 (IDEX) ---
$ws.Range("D40").Value = "'0.03883"

# --- Row 41 (KickToken) ---
$ws.Range("D41").Value = "'0.006694"

# --- Row 42 (BKEXToken) ---
$ws.Range("D42").Value = "'0.1071"

# --- Row 43 (CEJI) ---
$ws.Range("D43").Value = "'0.002900"

# --- Row 44 (LocalTraders) ---
$ws.Range("D44").Value = "'0.006502"

# --- Row 45 (CoinLion) ---
$ws.Range("D45").Value = "'0.00005608"

# --- Row 46 (Kangarootoken) ---
$ws.Range("D46").Value = "'0.00000000750"

# --- Row 47 (CoinbaseStockToken) ---
$ws.Range("D47").Value = "'0.3899"

# --- Row 49 (CryptobidCoin) ---
$ws.Range("D49").Value = "'0.00002100"

# --- Row 50 (SpecialPowerGold) ---
$ws.Range("D50").Value = "'0.01010"

# Cells that were given a leading text-quote above (because their new
# value looks numeric) get an implicit "Text" number format/style from
# the quote prefix. Put them back on the default "Normal" style so the
# workbook's styling is untouched, matching every other text cell here.
$numericTextCells = @( `
    "D2","D3","D4","D5","D7","D8","D9","D10","D11","D12","D13","D14", `
    "D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25", `
    "D26","D27","D40","D41","D42","D43","D44","D45","D46","D47","D49","D50" `
)
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).Style = "Normal"
}

Write-Host "Applied symbol-list update: $($numericTextCells.Count) price cells + 37 text cells."
